# Apply the per-cell price/volume (and one coin-name/link swap) updates
# recorded in the Thu Mar 23 21:25:52 UTC 2023 "Updated cryptos list" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.128.52"
$ws.Range("E2").Value = "  +2.66%  "
$ws.Range("D3").Value = "1.813.61"
$ws.Range("E3").Value = "  +4.40%  "
$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'329.15"
$ws.Range("E5").Value = "  +1.71%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").Value = "'0.4404"
$ws.Range("E7").Value = "  +3.96%  "
$ws.Range("D8").Value = "'0.3692"
$ws.Range("E8").Value = "  +2.18%  "
$ws.Range("D9").Value = "'44.66"
$ws.Range("E9").Value = "  -0.81%  "
$ws.Range("D10").Value = "'0.07681"
$ws.Range("E10").Value = "  +3.54%  "
$ws.Range("D11").Value = "'1.123"
$ws.Range("E11").Value = "  +0.12%  "
$ws.Range("D12").Value = "'1.005"
$ws.Range("E12").Value = "  +0.35%  "
$ws.Range("D13").Value = "'21.97"
$ws.Range("E13").Value = "  +1.77%  "
$ws.Range("D14").Value = "'6.250"
$ws.Range("E14").Value = "  +3.09%  "
$ws.Range("D15").Value = "'7.538"
$ws.Range("E15").Value = "  +5.30%  "
$ws.Range("D16").Value = "1.817.93"
$ws.Range("E16").Value = "  +4.60%  "
$ws.Range("D17").Value = "'92.46"
$ws.Range("E17").Value = "  +6.67%  "
$ws.Range("D18").Value = "'0.00001082"
$ws.Range("E18").Value = "  +1.64%  "
$ws.Range("D19").Value = "'0.06624"
$ws.Range("E19").Value = "  +10.20%  "
$ws.Range("D20").Value = "'1.002"
$ws.Range("E20").Value = "  +0.20%  "
$ws.Range("D21").Value = "'17.51"
$ws.Range("E21").Value = "  +4.13%  "
$ws.Range("D22").Value = "'6.196"
$ws.Range("E22").Value = "  +2.34%  "
$ws.Range("D23").Value = "28.230.13"
$ws.Range("E23").Value = "  +2.82%  "
$ws.Range("D24").Value = "'11.66"
$ws.Range("E24").Value = "  +3.30%  "
$ws.Range("D25").Value = "'2.045"
$ws.Range("E25").Value = "  -14.93%  "
$ws.Range("E26").Value = "  +2.84%  "
$ws.Range("D27").Value = "'155.66"
$ws.Range("E27").Value = "  +4.12%  "
$ws.Range("D28").Value = "2.023.14"
$ws.Range("E28").Value = "  +4.44%  "
$ws.Range("D29").Value = "'2.309"
$ws.Range("E29").Value = "  -2.09%  "
$ws.Range("D30").Value = "'128.27"
$ws.Range("E30").Value = "  +1.66%  "
$ws.Range("D31").Value = "'1.200"
$ws.Range("E31").Value = "  -5.40%  "
$ws.Range("D32").Value = "'5.854"
$ws.Range("E32").Value = "  +4.73%  "
$ws.Range("D33").Value = "'0.09210"
$ws.Range("E33").Value = "  +2.13%  "
$ws.Range("D34").Value = "'3.673"
$ws.Range("E34").Value = "  -1.43%  "
$ws.Range("D35").Value = "'13.02"
$ws.Range("E35").Value = "  +4.30%  "
$ws.Range("D36").Value = "'0.02346"
$ws.Range("E36").Value = "  +2.86%  "
$ws.Range("D37").Value = "'0.2168"
$ws.Range("E37").Value = "  +0.72%  "
$ws.Range("D38").Value = "'0.06208"
$ws.Range("E38").Value = "  +1.00%  "
$ws.Range("D39").Value = "'5.142"
$ws.Range("E39").Value = "  +2.42%  "
$ws.Range("D40").Value = "'0.6549"
$ws.Range("E40").Value = "  +2.25%  "
$ws.Range("D41").Value = "'1.196"
$ws.Range("E41").Value = "  +1.20%  "
$ws.Range("D42").Value = "'8.137"
$ws.Range("E42").Value = "  +3.95%  "
$ws.Range("D43").Value = "'1.001"
$ws.Range("E43").Value = "  +0.10%  "
$ws.Range("D44").Value = "'13.92"
$ws.Range("E44").Value = "  +2.70%  "
$ws.Range("D45").Value = "'1.387"
$ws.Range("E45").Value = "  -2.08%  "
$ws.Range("D46").Value = "'0.6063"
$ws.Range("E46").Value = "  +3.65%  "
$ws.Range("D47").Value = "'3.758"
$ws.Range("E47").Value = "  +0.34%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "'126.86"
$ws.Range("E48").Value = "  +1.11%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "'2.031"
$ws.Range("E49").Value = "  +5.05%  "
$ws.Range("D50").Value = "'1.152"
$ws.Range("E50").Value = "  +5.33%  "
$ws.Range("D51").Value = "'0.06985"
$ws.Range("E51").Value = "  +2.46%  "

# The quote-prefix from the apostrophe trick above leaves a "quote
# prefix" cell style; reset those cells back to the Normal style so
# only the text content (not formatting) changes.
$ws.Range("D4","D5","D6","D7","D8","D9","D10","D11","D12","D13","D14","D15","D17","D18","D19","D20","D21","D22","D24","D25","D27","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51").Style = "Normal"
